# Atualização de bases das ligas, do dia: 13-06-2024 às 19:35
# Swap the (mis-ordered) data rows so matches are listed in the correct order.
# Column A (the running index) stays put; columns B:AD (id .. PL_AhUnder) swap
# between each pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($sheet, $row1, $row2) {
    $range1 = $sheet.Range($sheet.Cells.Item($row1, 2), $sheet.Cells.Item($row1, 30))
    $range2 = $sheet.Range($sheet.Cells.Item($row2, 2), $sheet.Cells.Item($row2, 30))

    # NOTE: use Value2 for the read - Value's getter on a multi-cell Range
    # does not reliably marshal the array in this host.
    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value = $vals2
    $range2.Value = $vals1
}

Swap-Rows $ws 176 177
Swap-Rows $ws 187 188
Swap-Rows $ws 305 306
